$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Evidencias"

# Header row
$headers = @("Id","Status","Tipo","Quantidade","Notas","Profundidade","Solo","Id do Ponto","Longitude","Latitude","Status do Ponto","Usuário que criou","Usuário que atualizou","Data de criação","Data de atualização")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data rows
$data = @(
    @('"5dc8ab4637163e00178bd0e6"', "Em aberto", "Madeira", 1, "some fixed notes...", "Superfície", "Água", $null, -45.9776900762719, -20.7172729395149, "Realizado", '"5d16de7d8db2ea00174a916a"', '"5d16de7d8db2ea00174a916a"', 43775.04038826389, 43780.01979832176),
    @('"5dc8aae737163e00178bd0e4"', "Em aberto", "Madeira", 1, "some fixed notes...", "20-30cm", "Rocha", $null, -46.0381655465671, -20.6904377923282, "Realizado", '"5d16de7d8db2ea00174a916a"', '"5d16de7d8db2ea00174a916a"', 43775.04038831018, 43780.01895534722),
    @('"5dc3065937a84c0017f7bd1a"', "Em aberto", "Madeira", 1, "some fixed notes...", "Superfície", "Água", $null, -46.0478527038117, -20.6878642170664, "Realizado", '"5d16de7d8db2ea00174a916a"', '"5d16de7d8db2ea00174a916a"', 43775.040388333335, 43775.73864273148),
    @('"5dcaaf0d50f17900176dff76"', "Em aberto", "Madeira", 1, "some fixed notes...", "10-20cm", "Água", $null, -46.0308166248642, -20.6865044291569, "Realizado", '"5d16de7d8db2ea00174a916a"', '"5dc75f42ea7b0500177d4381"', 43775.0403883449, 43781.54820726852)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $value = $row[$c]
        if ($null -eq $value) {
            continue
        }
        $cell = $ws.Cells.Item($excelRow, $c + 1)
        $cell.Value = $value
        if ($c -eq 13 -or $c -eq 14) {
            # Data de criação / Data de atualização columns (N, O) -> date serials,
            # formatted with the built-in short-date number format (numFmtId 14).
            $cell.NumberFormat = "mm-dd-yy"
        }
    }
}

# Column widths: col 1 stays (OOXML) width 36, cols 2-15 (OOXML) width 50.
# Excel's ColumnWidth (character units) differs from the stored OOXML <col width>
# by a constant ~5/6 padding offset for the default Calibri 11 font, so we back
# that offset out to land on the exact target OOXML width.
$padding = 5 / 6
$ws.Columns.Item(1).ColumnWidth = 36 - $padding
for ($c = 2; $c -le 15; $c++) {
    $ws.Columns.Item($c).ColumnWidth = 50 - $padding
}
